# Zeitabrechnung, Sprint Report update
# Fills in Yichi Zhang's (column C) missing hours for early-December entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row -> date (for reference):
#  50 = 2021-12-04   55 = 2021-12-09   59 = 2021-12-13  60 = 2021-12-14
#  62 = 2021-12-16   63 = 2021-12-17   65 = 2021-12-19
$ws.Range("C50").Value = 1
$ws.Range("C55").Value = 1.5
$ws.Range("C59").Value = 3
$ws.Range("C60").Value = 4.5
$ws.Range("C62").Value = 4
$ws.Range("C63").Value = 1.5
$ws.Range("C65").Value = 9

# Drop the unused "Heading1" cell style (not referenced by any cell) and
# align the remaining auto-generated style names with their Excel defaults.
$wb.Styles.Item("Heading1").Delete()
$wb.Styles.Item("Result2").Name = "Ergebnis 2"
$wb.Styles.Item("Heading").Name = "Heading 3"

# Leave the selection where the edits finished.
$ws.Range("C66").Select()
